$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting of an existing date cell to the new date cells first
$ws.Range("B60").Copy()
$ws.Range("B61:B62").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 61
$ws.Range("A61").Value = "Federico Speroni"
$ws.Range("B61").Value = (Get-Date -Year 2017 -Month 5 -Day 30 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C61").Value = 3
$ws.Range("E61").Value = "Administrador"
$ws.Range("D61").Value = "Sprint 3 - BackEnd y FrontEnd"

# Row 62
$ws.Range("A62").Value = "Federico Speroni"
$ws.Range("B62").Value = (Get-Date -Year 2017 -Month 6 -Day 6 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C62").Value = 3
$ws.Range("E62").Value = "Administrador"
$ws.Range("D62").Value = "Sprint 3 - BackEnd y FrontEnd"

$ws.Range("D62").Select()
